$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix typo in row 17 description: "Mr.Chavais" -> "Mr.Chaveys"
$ws.Range("G17").Value = "Retour sur la fin du sprinte 1`net le début du sprinte 2 avec  Mr.Chaveys et Simon Cuany"

# 2. Fill in the previously-blank row 25 with new journal entry data
$ws.Range("B25").Value = 44277
$ws.Range("C25").Value = "10H40"
$ws.Range("D25").Value = "12h15"
$ws.Range("E25").Value = "95min"
$ws.Range("G25").Value = "Redaction du CDC "

# 3. Update the view state (scroll position + active cell/selection)
$ws.Application.ActiveWindow.ScrollRow = 14
$ws.Range("G25").Select()
